$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "Sheet2"

# Row 1: report generation date/time moved forward
$ws.Range("D1").Value = 45572
$ws.Range("F1").Value = 0.806519050925926

# Row 16/17: rename "Induction Hardening Bearing Surface 1" process to the
# plural "Induction Hardening Bearing Surfaces 1, 2" and update the impact
# figures for the less impactful hardening process
$ws.Range("U16").Value = "Induction Hardening Bearing Surfaces 1, 2"
$ws.Range("C17").Value = 174.399516254188
$ws.Range("U17").Value = 34.074702596165
